$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been recorded.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for both movies.`n"
$ws.Range("D3").Value = "both_movies, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded successfully.`n"
$ws.Range("D4").Value = "both_movies, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday's assembly has concluded without a selection.`n"
$ws.Range("D5").Value = "no_decision, "
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie has been selected for Friday.`n"
$ws.Range("D6").Value = "no_decision, "
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made.`n"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie.`"`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday was not made, indicating that no consensus was reached.`n"
$ws.Range("D9").Value = "no_decision, "
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie`" for the movie to be shown on Friday.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to be shown on Friday.`n"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected in this meeting.`n"
$ws.Range("D13").Value = "no_decision, "
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The function has been executed, indicating that no decision about Friday's movie was made.`n"
$ws.Range("D15").Value = "no_decision, "
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision concludes with no agreement on a movie for Friday, and no selection was made.`n"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for the movie `"Barbie`" for the Friday showing.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was officially selected for Friday.`n"
$ws.Range("D18").Value = "no_decision, "
$ws.Range("C19").Value = "MSG: None`n`nMSG: The committee has not reached a decision regarding which movie to show on Friday.`n"
$ws.Range("D19").Value = "no_decision, "
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision process has concluded without selecting a movie for Friday, and no film will be acquired.`n"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision about Friday’s movie can be made.`n"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has concluded without a selection.`n"
$ws.Range("D22").Value = "no_decision, "
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("D23").Value = "no_decision, "
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been made.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached, resulting in no selection being made.`n"
$ws.Range("D25").Value = "no_decision, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for showing on Friday.`n"
$ws.Range("D26").Value = "both_movies, "
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has been recorded as `"no decision.`"`n"
$ws.Range("D27").Value = "no_decision, "
$ws.Range("C28").Value = "MSG: None`n`nMSG: The committee did not arrive at a decision regarding which movie to show on Friday.`n"
$ws.Range("D28").Value = "no_decision, "
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be shown on Friday.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday has not been made.`n"
$ws.Range("D31").Value = "no_decision, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("D32").Value = "no_decision, "
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The conversation did not lead to a decision about which movie to show on Friday. Therefore, the appropriate action is to acknowledge that no decision can be made at this time.`n"
$ws.Range("D34").Value = "no_decision, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D35").Value = "no_decision, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded, and there is no consensus on which movie to show on Friday.`n"
$ws.Range("D37").Value = "no_decision, "
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been noted, and no movie has been selected for Friday's screening.`n"
$ws.Range("D38").Value = "no_decision, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D39").Value = "both_movies, "
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded: `"Barbie`" was successfully selected for the screening.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made, resulting in no decision overall.`n"
$ws.Range("D41").Value = "no_decision, "
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("D42").Value = "no_decision, "
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has resulted in no conclusive agreement.`n"
$ws.Range("D43").Value = "no_decision, "
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision-making committee did not reach a consensus regarding which movie to show on Friday, resulting in no decision being made.`n"
$ws.Range("D44").Value = "no_decision, "
$ws.Range("C45").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("D45").Value = "both_movies, "
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: No decision has been made regarding the movie to be shown on Friday.`n"
$ws.Range("D48").Value = "no_decision, "
$ws.Range("C49").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has not been reached, so I will call the no_decision function.`n"
$ws.Range("D50").Value = "no_decision, "
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie.`"`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected in the meeting.`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision made regarding the movie for Friday.`n"
$ws.Range("D53").Value = "no_decision, "
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("D54").Value = "no_decision, "
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision-making process did not yield a choice for Friday's movie, resulting in no decision being made.`n"
$ws.Range("D55").Value = "no_decision, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision process has concluded without a decision on which movie to acquire for Friday.`n"
$ws.Range("D56").Value = "no_decision, "
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision has been recorded: `"Barbie`" was selected for the movie to be shown on Friday.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to be shown on Friday.`n"
$ws.Range("D58").Value = "no_decision, "
$ws.Range("C59").Value = "MSG: None`n`nMSG: The conversation ended without a decision on which movie to show on Friday.`n"
$ws.Range("D59").Value = "no_decision, "
$ws.Range("C60").Value = "MSG: None`n`nMSG: The conversation resulted in no decision about which movie will be shown on Friday.`n"
$ws.Range("D60").Value = "no_decision, "
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday resulted in no agreement.`n"
$ws.Range("D61").Value = "no_decision, "
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded.`n"
$ws.Range("D64").Value = "both_movies, "
$ws.Range("C65").Value = "MSG: None`n`nMSG: The function has been called to indicate that no decision was made regarding the movie selection.`n"
$ws.Range("D65").Value = "no_decision, "
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision has been made that no movie will be shown on Friday.`n"
$ws.Range("D66").Value = "no_decision, "
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding Friday's movie.`n"
$ws.Range("D67").Value = "no_decision, "
$ws.Range("C68").Value = "MSG: None`n`nMSG: The decision has been recorded indicating that no movie was selected during the discussion.`n"
$ws.Range("D68").Value = "no_decision, "

Write-Host "Applied judgement log updates"
